$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '26.979.65'
Set-TextValue 'E2' '  +0.18%  '
Set-TextValue 'D3' '1.681.97'
Set-TextValue 'E3' '  +0.68%  '
Set-TextValue 'E4' '  -0.13%  '
Set-TextValue 'D5' '215.22'
Set-TextValue 'E5' '  -0.30%  '
Set-TextValue 'E6' '  -2.33%  '
Set-TextValue 'E7' '  -0.17%  '
Set-TextValue 'E8' '  -1.01%  '
Set-TextValue 'E9' '  -0.03%  '
Set-TextValue 'D10' '20.97'
Set-TextValue 'E10' '  +3.89%  '
Set-TextValue 'D11' '0.0887'
Set-TextValue 'E11' '  -0.30%  '
Set-TextValue 'D12' '1.917.32'
Set-TextValue 'E12' '  +0.55%  '
Set-TextValue 'D13' '1.677.68'
Set-TextValue 'E13' '  +0.48%  '
Set-TextValue 'D14' '4.11'
Set-TextValue 'E14' '  +0.47%  '
Set-TextValue 'E15' '  +2.25%  '
Set-TextValue 'D16' '65.89'
Set-TextValue 'E16' '  +0.46%  '
Set-TextValue 'B17' 'Chainlink'
Set-TextValue 'C17' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D17' '8.21'
Set-TextValue 'E17' '  +5.71%  '
Set-TextValue 'B18' 'WrappedBTC'
Set-TextValue 'C18' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue 'D18' '27.000.65'
Set-TextValue 'E18' '  +0.17%  '
Set-TextValue 'D19' '236.23'
Set-TextValue 'E19' '  +1.70%  '
Set-TextValue 'D20' '0.0₃0734'
Set-TextValue 'E20' '  +0.04%  '
Set-TextValue 'E21' '  -0.14%  '
Set-TextValue 'E22' '  -0.11%  '
Set-TextValue 'D23' '9.22'
Set-TextValue 'E23' '  +0.33%  '
Set-TextValue 'E24' '  -3.97%  '
Set-TextValue 'D25' '146.42'
Set-TextValue 'E25' '  +0.56%  '
Set-TextValue 'E26' '  +1.22%  '
Set-TextValue 'D27' '16.09'
Set-TextValue 'E27' '  +1.53%  '
Set-TextValue 'E28' '  -2.98%  '
Set-TextValue 'E29' '  +0.04%  '
Set-TextValue 'E30' '  +0.68%  '
Set-TextValue 'E31' '  -0.55%  '
Set-TextValue 'E32' '  +0.82%  '
Set-TextValue 'D33' '1.494.52'
Set-TextValue 'E33' '  +2.85%  '
Set-TextValue 'E34' '  +1.12%  '
Set-TextValue 'E35' '  +4.90%  '
Set-TextValue 'E36' '  -0.12%  '
Set-TextValue 'D37' '0.586'
Set-TextValue 'E37' '  +3.73%  '
Set-TextValue 'D38' '0.917'
Set-TextValue 'E38' '  +2.19%  '
Set-TextValue 'E39' '  +3.88%  '
Set-TextValue 'E40' '  +6.65%  '
Set-TextValue 'E41' '  -4.94%  '
Set-TextValue 'E42' '  -0.13%  '
Set-TextValue 'D43' '67.56'
Set-TextValue 'E43' '  +2.85%  '
Set-TextValue 'D44' '2.28'
Set-TextValue 'E44' '  -1.38%  '
Set-TextValue 'D45' '1.821.97'
Set-TextValue 'E45' '  +0.71%  '
Set-TextValue 'D46' '0.780'
Set-TextValue 'E46' '  +0.36%  '
Set-TextValue 'D47' '90.66'
Set-TextValue 'E47' '  +0.04%  '
Set-TextValue 'D48' '0.0₆0104'
Set-TextValue 'E48' '  -0.69%  '
Set-TextValue 'D49' '1.53'
Set-TextValue 'E49' '  -0.37%  '
Set-TextValue 'E50' '  +3.36%  '
Set-TextValue 'B51' 'Cronos'
Set-TextValue 'C51' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D51' '0.0508'
Set-TextValue 'E51' '  +0.03%  '
